# Add a quadrilateral ("Triang") element's displacement results alongside
# the existing ones: relabel the header row so the X/Y/UX/UY columns read
# as the quad-element variants, and widen columns B & C to fit the longer
# headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = "X (Triang)"
$ws.Range("C1").Value = "Y (Triang)"
$ws.Range("D1").Value = "UX (Triang)"
$ws.Range("E1").Value = "UY (Triang)"

# Widen columns B and C so the longer headers fit (closest the engine's
# character->pixel quantization allows to the target 9.57 / 9.43 widths).
$ws.Columns.Item(2).ColumnWidth = 8.6666666666667
$ws.Columns.Item(3).ColumnWidth = 8.6666666666667
